$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text (non-numeric-looking) cell updates: Coin names, URLs, Volume codes ---
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("E10").Value = '9WazirXWRX'
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("E11").Value = '10MandalaExchangeTokenMDX'
$ws.Range("B12").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C12").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("E12").Value = '11LiechtensteinCryptoassetsExchangeLCX'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("E13").Value = '12BitrueCoinBTR'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("E14").Value = '13BitMartTokenBMX'
$ws.Range("B15").Value = 'MCDex'
$ws.Range("C15").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("E15").Value = '14MCDexMCB'
$ws.Range("B16").Value = 'BitForexToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("E16").Value = '15BitForexTokenBF'
$ws.Range("B17").Value = 'CoinExToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("E17").Value = '16CoinExTokenCET'
$ws.Range("B18").Value = 'UpBots'
$ws.Range("C18").Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range("E18").Value = '17UpBotsUBXTBestin24h'
$ws.Range("B19").Value = 'One'
$ws.Range("C19").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("E19").Value = '18OneONEWorstin24h'
$ws.Range("E47").Value = '46CoinbaseStockTokenCOIN'

# --- Numeric-looking Price cells: force text format so exact string (incl. trailing zeros) is preserved ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '249.04'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '22.64'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.258'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05689'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '3.408'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '6.342'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8057'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9058'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1401'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07443'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03103'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03010'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09382'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.867'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.001570'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.04774'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.01828'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0005801'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.006452'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.004994'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0009990'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.697'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.195'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.3260'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03983'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006711'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1070'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002766'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.007721'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005594'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.2060'
